$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 79.319999999999993
$ws.Range("F2").Value = -2.0740740740740824
$ws.Range("G2").Value = $false

$ws.Range("C3").Value = 9792.59
